$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for the D-column price cells whose new values would
# otherwise be auto-recognised by Excel as plain numbers (they contain a
# single decimal point), so they stay text just like the rest of that
# column (e.g. "69.477.59" is untouched because Excel can't parse it as a
# number anyway).
$textRows = @(5,6,10,14,18,19,20,22,23,28,32,33,34,37,38,39,40,42,45,46,47,49,50)
foreach ($r in $textRows) {
    $ws.Range("D$r").NumberFormat = "@"
}

$ws.Range("D2").Value = "69.477.59"
$ws.Range("E2").Value = "  +3.15%  "
$ws.Range("D3").Value = "3.372.62"
$ws.Range("E3").Value = "  +4.57%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "191.79"
$ws.Range("E5").Value = "  +4.39%  "
$ws.Range("D6").Value = "592.73"
$ws.Range("E6").Value = "  +2.49%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("E9").Value = "  +2.97%  "
$ws.Range("D10").Value = "6.77"
$ws.Range("E10").Value = "  +3.19%  "
$ws.Range("E11").Value = "  +2.64%  "
$ws.Range("D12").Value = "3.960.53"
$ws.Range("E12").Value = "  +4.77%  "
$ws.Range("E13").Value = "  +1.39%  "
$ws.Range("D14").Value = "28.65"
$ws.Range("E14").Value = "  +3.60%  "
$ws.Range("D15").Value = "69.509.09"
$ws.Range("E15").Value = "  +3.10%  "
$ws.Range("E16").Value = "  +2.05%  "
$ws.Range("D17").Value = "3.385.05"
$ws.Range("E17").Value = "  +5.04%  "
$ws.Range("D18").Value = "450.86"
$ws.Range("E18").Value = "  +14.06%  "
$ws.Range("D19").Value = "5.85"
$ws.Range("E19").Value = "  +1.72%  "
$ws.Range("D20").Value = "13.83"
$ws.Range("E20").Value = "  +3.01%  "
$ws.Range("E21").Value = "  +3.70%  "
$ws.Range("D22").Value = "74.54"
$ws.Range("E22").Value = "  +4.86%  "
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("D24").Value = "3.517.84"
$ws.Range("E25").Value = "  +4.78%  "
$ws.Range("E26").Value = "  +0.80%  "
$ws.Range("E27").Value = "  +3.60%  "
$ws.Range("D28").Value = "9.53"
$ws.Range("E28").Value = "  +0.31%  "
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("E30").Value = "  +2.18%  "
$ws.Range("E31").Value = "  +3.32%  "
$ws.Range("D32").Value = "5.66"
$ws.Range("E32").Value = "  +2.20%  "
$ws.Range("D33").Value = "1.29"
$ws.Range("E33").Value = "  +3.04%  "
$ws.Range("D34").Value = "7.01"
$ws.Range("E34").Value = "  +0.78%  "
$ws.Range("E36").Value = "  +5.03%  "
$ws.Range("D37").Value = "165.63"
$ws.Range("E37").Value = "  +3.10%  "
$ws.Range("D38").Value = "1.95"
$ws.Range("E38").Value = "  +3.42%  "
$ws.Range("D39").Value = "27.33"
$ws.Range("E39").Value = "  +3.61%  "
$ws.Range("D40").Value = "0.819"
$ws.Range("E40").Value = "  +2.18%  "
$ws.Range("E41").Value = "  +1.80%  "
$ws.Range("D42").Value = "6.55"
$ws.Range("E42").Value = "  +0.74%  "
$ws.Range("D43").Value = "2.739.67"
$ws.Range("E43").Value = "  +5.60%  "
$ws.Range("E44").Value = "  +3.66%  "
$ws.Range("D45").Value = "25.71"
$ws.Range("E45").Value = "  +4.96%  "
$ws.Range("D46").Value = "0.0691"
$ws.Range("E46").Value = "  +1.23%  "
$ws.Range("D47").Value = "344.34"
$ws.Range("E47").Value = "  +3.08%  "
$ws.Range("E48").Value = "  +0.52%  "
$ws.Range("D49").Value = "0.0286"
$ws.Range("E49").Value = "  +3.17%  "
$ws.Range("D50").Value = "32.94"
$ws.Range("E50").Value = "  +8.10%  "
$ws.Range("E51").Value = "  +5.80%  "
